$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is being inserted above the existing row 23
# (Vega Central Mapocho de Santiago - Orégano, 2020-12-16 entry).
# Copy row 23 (so formatting / style carries over, including the
# date number format on column D) and insert it above itself; this
# shifts the old rows 23-45 down to 24-46 and bumps the sheet's used
# range to A1:R46.
$ws.Rows("23").Copy()
$ws.Rows("23").Insert()

# Now overwrite the new row 23 with the actual new observation's
# date (2022-01-05) and volume; the remaining columns (K/L/M price
# stats, etc.) match what was already copied from the old row 23.
$ws.Cells.Item(23, 4).Value = 44566
$ws.Cells.Item(23, 10).Value = 16
